$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 81.01587301587301
$ws.Cells.Item(3, 7).Value = 81.01587301587301
$ws.Cells.Item(4, 7).Value = 81.01587301587301
$ws.Cells.Item(5, 7).Value = 81.01587301587301
$ws.Cells.Item(6, 7).Value = 81.01587301587301
$ws.Cells.Item(7, 7).Value = 81.01587301587301
$ws.Cells.Item(8, 7).Value = 81.01587301587301
$ws.Cells.Item(9, 7).Value = 81.01587301587301
$ws.Cells.Item(10, 7).Value = 81.01587301587301
$ws.Cells.Item(11, 4).Value = 0.5174124598577119
$ws.Cells.Item(11, 7).Value = 81
$ws.Cells.Item(12, 4).Value = 0.5161807417094764
$ws.Cells.Item(12, 7).Value = 80.96825396825396
$ws.Cells.Item(13, 4).Value = 0.5161807417094764
$ws.Cells.Item(13, 7).Value = 80.96825396825396
$ws.Cells.Item(14, 4).Value = 0.5161807417094764
$ws.Cells.Item(14, 7).Value = 80.96825396825396
$ws.Cells.Item(15, 4).Value = 0.5161807417094764
$ws.Cells.Item(15, 7).Value = 80.96825396825396
$ws.Cells.Item(16, 4).Value = 0.5161807417094764
$ws.Cells.Item(16, 7).Value = 80.96825396825396
$ws.Cells.Item(17, 4).Value = 0.5161807417094764
$ws.Cells.Item(17, 7).Value = 80.96825396825396
$ws.Cells.Item(18, 4).Value = 0.5152274224425172
$ws.Cells.Item(18, 7).Value = 80.92063492063492
$ws.Cells.Item(19, 4).Value = 0.5152274224425172
$ws.Cells.Item(19, 7).Value = 80.92063492063492
$ws.Cells.Item(20, 4).Value = 0.5141847639893176
$ws.Cells.Item(20, 7).Value = 80.85714285714286
$ws.Cells.Item(21, 4).Value = 0.5136714286906973
$ws.Cells.Item(21, 7).Value = 80.80952380952381
$ws.Cells.Item(22, 4).Value = 0.5135967960445158
$ws.Cells.Item(22, 7).Value = 80.76190476190476
$ws.Cells.Item(23, 4).Value = 0.513614710902777
$ws.Cells.Item(23, 7).Value = 80.73015873015873
$ws.Cells.Item(24, 4).Value = 0.5134998964503081
$ws.Cells.Item(24, 7).Value = 80.68253968253968
$ws.Cells.Item(25, 4).Value = 0.5129417063919948
$ws.Cells.Item(25, 7).Value = 80.58730158730158
$ws.Cells.Item(26, 4).Value = 0.5121013756995278
$ws.Cells.Item(26, 7).Value = 80.49206349206349
$ws.Cells.Item(27, 4).Value = 0.5118789027780003
$ws.Cells.Item(27, 7).Value = 80.31746031746032
$ws.Cells.Item(28, 4).Value = 0.5091293873233976
$ws.Cells.Item(28, 7).Value = 80.11111111111111
$ws.Cells.Item(29, 4).Value = 0.5088898849734458
$ws.Cells.Item(29, 7).Value = 80.04761904761905
$ws.Cells.Item(30, 4).Value = 0.5076022467249167
$ws.Cells.Item(30, 7).Value = 79.7936507936508
$ws.Cells.Item(31, 4).Value = 0.5069011502526382
$ws.Cells.Item(31, 7).Value = 79.50793650793651
$ws.Cells.Item(32, 4).Value = 0.5062890186272081
$ws.Cells.Item(32, 7).Value = 79.03174603174604
$ws.Cells.Item(33, 4).Value = 0.504221301073956
$ws.Cells.Item(33, 7).Value = 78.58730158730158
$ws.Cells.Item(34, 4).Value = 0.5036122993922157
$ws.Cells.Item(34, 7).Value = 78.17460317460318
$ws.Cells.Item(35, 4).Value = 0.5016961756989567
$ws.Cells.Item(35, 7).Value = 77.58730158730158
$ws.Cells.Item(36, 4).Value = 0.5001333506558113
$ws.Cells.Item(36, 7).Value = 76.98412698412699
$ws.Cells.Item(37, 4).Value = 0.4987946814064388
$ws.Cells.Item(37, 7).Value = 76.22222222222223
$ws.Cells.Item(38, 4).Value = 0.4974823021793375
$ws.Cells.Item(38, 7).Value = 75.49206349206349
$ws.Cells.Item(39, 4).Value = 0.4962179496825814
$ws.Cells.Item(39, 7).Value = 74.47619047619048
$ws.Cells.Item(40, 4).Value = 0.4968859731646126
$ws.Cells.Item(40, 7).Value = 73.68253968253968
$ws.Cells.Item(41, 4).Value = 0.4959461033568031
$ws.Cells.Item(41, 7).Value = 72.3015873015873
$ws.Cells.Item(42, 4).Value = 0.4964571978195787
$ws.Cells.Item(42, 7).Value = 71.07936507936508
$ws.Cells.Item(43, 4).Value = 0.497199142789752
$ws.Cells.Item(43, 7).Value = 69.63492063492063
$ws.Cells.Item(44, 4).Value = 0.4981281830615387
$ws.Cells.Item(44, 7).Value = 67.87301587301587
$ws.Cells.Item(45, 4).Value = 0.5003782309121513
$ws.Cells.Item(45, 7).Value = 65.63492063492063
$ws.Cells.Item(46, 4).Value = 0.4666810476151583
$ws.Cells.Item(46, 7).Value = 60.79365079365079
$ws.Cells.Item(47, 4).Value = 0.3743476107988064
$ws.Cells.Item(47, 7).Value = 54.61904761904762
$ws.Cells.Item(48, 4).Value = 0.2904281899469116
$ws.Cells.Item(48, 7).Value = 47.31746031746032
$ws.Cells.Item(49, 4).Value = 0.241605468697636
$ws.Cells.Item(49, 7).Value = 32.25396825396825
$ws.Cells.Item(50, 4).Value = 0.2491353158916866
$ws.Cells.Item(50, 7).Value = 10.19047619047619

# Clear D51 (rmse_z_ga) - no longer has a computed value for this row
$ws.Cells.Item(51, 4).ClearContents()

